$d = $word.ActiveDocument

# 1. Merge "Meus em" + "préstimos" -> "Meus empréstimos"
$d.Content.Find.Execute("Meus empréstimos", $true, $false, $false, $false, $false, $true, 1, $false, "Meus empréstimos", 2)

# 2. Merge "Identificar caso haja empréstimos vencidos" + "/perto de vencer." into one run
$d.Content.Find.Execute("Identificar caso haja empréstimos vencidos/perto de vencer.", $true, $false, $false, $false, $false, $true, 1, $false, "Identificar caso haja empréstimos vencidos/perto de vencer.", 2)

# 3. Remove "2.6. " before "(A)." in the list item text
$d.Content.Find.Execute("como na figura 2.6. (A).", $true, $false, $false, $false, $false, $true, 1, $false, "como na figura (A).", 2)

# 4. Remove "2.6. " before "(B)." in the list item text
$d.Content.Find.Execute("como na figura 2.6. (B).", $true, $false, $false, $false, $false, $true, 1, $false, "como na figura (B).", 2)

# 5. Merge "A tela de Detalhes deve disponibilizar a infor" + "mações sobre o empréstimo como disposto no protótipo." into one run
$d.Content.Find.Execute("A tela de Detalhes deve disponibilizar a informações sobre o empréstimo como disposto no protótipo.", $true, $false, $false, $false, $false, $true, 1, $false, "A tela de Detalhes deve disponibilizar a informações sobre o empréstimo como disposto no protótipo.", 2)

# 6. Merge "A tela de Detalhes deve disponibilizar o botão “Voltar”, que per" + "mite voltar para a tela anterior." into one run
$d.Content.Find.Execute("A tela de Detalhes deve disponibilizar o botão “Voltar”, que permite voltar para a tela anterior.", $true, $false, $false, $false, $false, $true, 1, $false, "A tela de Detalhes deve disponibilizar o botão “Voltar”, que permite voltar para a tela anterior.", 2)

# 7. Merge "Figura " + "(A): " into one run (both bold)
$d.Content.Find.Execute("Figura (A): ", $true, $false, $false, $false, $false, $true, 1, $false, "Figura (A): ", 2)

# 8. Merge "Figura" + " (B): " into one run (both bold)
$d.Content.Find.Execute("Figura (B): ", $true, $false, $false, $false, $false, $true, 1, $false, "Figura (B): ", 2)
